# analysis_draft2017b: updating analysis to use marc_s2 simulation data
# Updates the sulfate-aerosol budget table (sources/sinks/burden/lifetime,
# per size bin) on Sheet1 with the refreshed marc_s2 simulation numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "+538.70 ± 0.75"
$ws.Range("D2").Value = "+0.09 ± 0.00"
$ws.Range("E2").Value = "+509.28 ± 0.75"
$ws.Range("F2").Value = "+26.45 ± 0.04"
$ws.Range("G2").Value = "+6.46 ± 0.03"
$ws.Range("B4").Value = "+11.55 ± 0.05"
$ws.Range("D4").Value = "+0.08 ± 0.00"
$ws.Range("E4").Value = "+4.37 ± 0.01"
$ws.Range("F4").Value = "+3.84 ± 0.02"
$ws.Range("G4").Value = "+3.25 ± 0.02"
$ws.Range("B5").Value = "+22.33 ± 0.02"
$ws.Range("F5").Value = "+20.59 ± 0.02"
$ws.Range("G5").Value = "+1.75 ± 0.00"
$ws.Range("E6").Value = "+0.09 ± 0.00"
$ws.Range("F7").Value = "+2.02 ± 0.01"
$ws.Range("G7").Value = "+1.47 ± 0.01"
$ws.Range("B8").Value = "+504.82 ± 0.75"
$ws.Range("E8").Value = "+504.82 ± 0.75"
$ws.Range("B9").Value = "-538.88 ± 0.74"
$ws.Range("D9").Value = "-0.09 ± 0.00"
$ws.Range("E9").Value = "-509.37 ± 0.75"
$ws.Range("F9").Value = "-26.52 ± 0.04"
$ws.Range("G9").Value = "-6.48 ± 0.03"
$ws.Range("D10").Value = "-0.09 ± 0.00"
$ws.Range("E11").Value = "-3.49 ± 0.02"
$ws.Range("B12").Value = "-396.60 ± 0.69"
$ws.Range("E12").Value = "-379.23 ± 0.68"
$ws.Range("F12").Value = "-14.19 ± 0.02"
$ws.Range("G12").Value = "-3.17 ± 0.01"
$ws.Range("B13").Value = "-20.75 ± 0.06"
$ws.Range("E13").Value = "-20.08 ± 0.06"
$ws.Range("F13").Value = "-0.56 ± 0.00"
$ws.Range("G13").Value = "-0.11 ± 0.00"
$ws.Range("B14").Value = "-116.29 ± 0.06"
$ws.Range("E14").Value = "-102.21 ± 0.06"
$ws.Range("F14").Value = "-11.08 ± 0.02"
$ws.Range("G14").Value = "-3.00 ± 0.01"
$ws.Range("B15").Value = "-5.25 ± 0.01"
$ws.Range("E15").Value = "-4.34 ± 0.01"
$ws.Range("F15").Value = "-0.70 ± 0.00"
$ws.Range("B16").Value = "+1.33 ± 0.00"
$ws.Range("F16").Value = "+0.32 ± 0.00"
$ws.Range("C17").Value = "+0.13 ± 0.00"
$ws.Range("D17").Value = "+0.09 ± 0.00"
$ws.Range("E17").Value = "+0.67 ± 0.00"
$ws.Range("F17").Value = "+4.39 ± 0.02"
$ws.Range("G17").Value = "+3.99 ± 0.04"
